$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers used by this "Generate Report for handoff" run.
# ---------------------------------------------------------------------------
$newMdName     = "7e073687-e8e9-4b70-aa2e-2cd31a43bdc0.md"
$newStatus     = "Handoff transform failed"
$epoch         = "0001-01-01 00:00:00"
$ignored       = "Ignored"
$configName    = ".localization-config"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/594bf3bb7715d4f5073f7f85da13736a254fc1b3/e2e/5895ebb2-82a2-4f45-aa04-ed8c3bec51ef.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/594bf3bb7715d4f5073f7f85da13736a254fc1b3/.localization-config"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c42f1c7857a60620b4c1758756625b0afa7253fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5895ebb2-82a2-4f45-aa04-ed8c3bec51ef.40f00867a549378d7206cf6d553519693767aeae.zh-cn.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cc3b4f4764e91f5998433ecb2ab2df57db2b6e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5895ebb2-82a2-4f45-aa04-ed8c3bec51ef.40f00867a549378d7206cf6d553519693767aeae.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet: rename the handed-off file and flip its status in both
# locale columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", $configName)

# ---------------------------------------------------------------------------
# Per-locale detail sheets: clear the stale handoff artifact, reset the
# handoff bookkeeping columns, and mark the row ignored.
# ---------------------------------------------------------------------------
function Update-LocaleSheet($ws, $xlfUrl) {
    $ws.Range("A2").Value = $newMdName
    $ws.Range("B2").Value = $newStatus
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = $epoch
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $ignored

    $ws.Range("D3").Value = $epoch
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $ignored

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $newMdName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $xlfUrl, "", "", $configName)
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn $zhXlfUrl

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe $deXlfUrl
